$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 85326300.61677904
$ws.Range("C2").Value = 99068887.23512885
$ws.Range("D2").Value = 112811473.8534783
$ws.Range("E2").Value = 126554060.47182785
$ws.Range("F2").Value = 140296647.09017748
$ws.Range("B3").Value = 198216755.83586156
$ws.Range("C3").Value = 211959342.45421135
$ws.Range("D3").Value = 225701929.07256082
$ws.Range("E3").Value = 239444515.69091037
$ws.Range("F3").Value = 253187102.30926
$ws.Range("B4").Value = 424234826.04219234
$ws.Range("C4").Value = 437977412.66054213
$ws.Range("D4").Value = 451719999.27889156
$ws.Range("E4").Value = 465462585.8972412
$ws.Range("F4").Value = 479205172.51559085
$ws.Range("B5").Value = 695874424.4650055
$ws.Range("C5").Value = 709617011.0833554
$ws.Range("D5").Value = 723359597.7017049
$ws.Range("E5").Value = 737102184.3200544
$ws.Range("F5").Value = 750844770.938404
